$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9083871841430664
$ws.Range("B1").Value = 1.966000199317932
$ws.Range("C1").Value = 3.035072326660156
$ws.Range("D1").Value = 3.695548534393311
$ws.Range("E1").Value = 1.761780142784119
